$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2999.1667
$ws.Range("I62").Value = 2998.3333
$ws.Range("K62").Value = 2998.3333
$ws.Range("M62").Value = -2374.3333

$ws.Range("H65").Value = 2999.1667
$ws.Range("I65").Value = 2998.3333
$ws.Range("K65").Value = 14991.6665
$ws.Range("M65").Value = -11871.6665

$ws.Range("H98").Value = 2194.9048
$ws.Range("I98").Value = 2331.842
$ws.Range("J98").Value = 894
$ws.Range("K98").Value = 2331.842
$ws.Range("L98").Value = 894
$ws.Range("M98").Value = -833.8420000000001
$ws.Range("N98").Value = -3890

$ws.Range("H118").Value = 495
$ws.Range("I118").Value = 495
$ws.Range("K118").Value = 1485
$ws.Range("M118").Value = 172

$ws.Range("H122").Value = 2194.9048
$ws.Range("I122").Value = 2331.842
$ws.Range("J122").Value = 894
$ws.Range("K122").Value = 6995.526
$ws.Range("L122").Value = 2682
$ws.Range("M122").Value = -4545.526
$ws.Range("N122").Value = -7582

$ws.Range("H129").Value = 860.1539
$ws.Range("J129").Value = 881.2826
$ws.Range("L129").Value = 2643.8478
$ws.Range("N129").Value = -12643.8478

$ws.Range("H132").Value = 1158.4048
$ws.Range("I132").Value = 980.2105
$ws.Range("J132").Value = 2851.25
$ws.Range("K132").Value = 2940.6315
$ws.Range("L132").Value = 8553.75
$ws.Range("M132").Value = -410.6315
$ws.Range("N132").Value = -13613.75

$ws.Range("H137").Value = 41584.12
$ws.Range("I137").Value = 1494.7646
$ws.Range("J137").Value = 126774
$ws.Range("K137").Value = 4484.293799999999
$ws.Range("L137").Value = 380322
$ws.Range("M137").Value = -1934.293799999999
$ws.Range("N137").Value = -385422

$ws.Range("H141").Value = 936005.7
$ws.Range("I141").Value = 1219153.5
$ws.Range("K141").Value = 3657460.5
$ws.Range("M141").Value = -3652280.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3913.5
$ws.Range("I61").Value = 2042.8572
$ws.Range("J61").Value = 6532.4
$ws.Range("K61").Value = 2042.8572
$ws.Range("L61").Value = 6532.4
$ws.Range("M61").Value = -1830.8572
$ws.Range("N61").Value = -6956.4

$ws.Range("H63").Value = 5491.8335
$ws.Range("J63").Value = 4149.8335
$ws.Range("L63").Value = 4149.8335
$ws.Range("N63").Value = -5521.8335

$ws.Range("H66").Value = 5491.8335
$ws.Range("J66").Value = 4149.8335
$ws.Range("L66").Value = 20749.1675
$ws.Range("N66").Value = -27613.1675

$ws.Range("H74").Value = 1345.0312
$ws.Range("I74").Value = 1130.3334
$ws.Range("J74").Value = 1989.125
$ws.Range("K74").Value = 1130.3334
$ws.Range("L74").Value = 1989.125
$ws.Range("M74").Value = -256.3334
$ws.Range("N74").Value = -3737.125

$ws.Range("H77").Value = 1345.0312
$ws.Range("I77").Value = 1130.3334
$ws.Range("J77").Value = 1989.125
$ws.Range("K77").Value = 5651.666999999999
$ws.Range("L77").Value = 9945.625
$ws.Range("M77").Value = -1283.666999999999
$ws.Range("N77").Value = -18681.625

$ws.Range("H132").Value = 1505.6765
$ws.Range("I132").Value = 1012.6667
$ws.Range("K132").Value = 3038.0001
$ws.Range("M132").Value = -508.0001000000002

$ws.Range("H136").Value = 3913.5
$ws.Range("I136").Value = 2042.8572
$ws.Range("J136").Value = 6532.4
$ws.Range("K136").Value = 6128.571599999999
$ws.Range("L136").Value = 19597.2
$ws.Range("M136").Value = -3578.571599999999
$ws.Range("N136").Value = -24697.2

$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 8616.65
$ws.Range("I134").Value = 10854.6
$ws.Range("K134").Value = 32563.8
$ws.Range("M134").Value = -30028.8

$ws.Range("H140").Value = 53837.5
$ws.Range("J140").Value = 53837.5
$ws.Range("L140").Value = 53837.5
$ws.Range("N140").Value = -64197.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2103.4614
$ws.Range("I31").Value = 1607
$ws.Range("K31").Value = 1607
$ws.Range("M31").Value = -1312

$ws.Range("H34").Value = 2103.4614
$ws.Range("I34").Value = 1607
$ws.Range("K34").Value = 1607
$ws.Range("M34").Value = -1405

$ws.Range("H58").Value = 1674271.2
$ws.Range("I58").Value = 2416940.5
$ws.Range("J58").Value = 3265.125
$ws.Range("K58").Value = 2416940.5
$ws.Range("L58").Value = 3265.125
$ws.Range("M58").Value = -2416737.5
$ws.Range("N58").Value = -3671.125

$ws.Range("H132").Value = 1844.1702
$ws.Range("I132").Value = 1186.7878
$ws.Range("K132").Value = 3560.3634
$ws.Range("M132").Value = -1030.3634

$ws.Range("H134").Value = 1981.0769
$ws.Range("I134").Value = 1717.7826
$ws.Range("K134").Value = 5153.3478
$ws.Range("M134").Value = -2618.3478

$ws.Range("H136").Value = 1674271.2
$ws.Range("I136").Value = 2416940.5
$ws.Range("J136").Value = 3265.125
$ws.Range("K136").Value = 7250821.5
$ws.Range("L136").Value = 9795.375
$ws.Range("M136").Value = -7248271.5
$ws.Range("N136").Value = -14895.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 167.875
$ws.Range("I2").Value = 198.8
$ws.Range("J2").Value = 116.333336
$ws.Range("K2").Value = 1192.8
$ws.Range("L2").Value = 698.000016
$ws.Range("M2").Value = -1079.8
$ws.Range("N2").Value = -924.000016

$ws.Range("H4").Value = 344903.34
$ws.Range("I4").Value = 344903.34
$ws.Range("K4").Value = 1034710.02
$ws.Range("M4").Value = -1034598.02

$ws.Range("H103").Value = 2410.5334
$ws.Range("I103").Value = 2301
$ws.Range("K103").Value = 6903
$ws.Range("M103").Value = -6024

$ws.Range("H117").Value = 15873434
$ws.Range("I117").Value = 432
$ws.Range("J117").Value = 28571836
$ws.Range("K117").Value = 1296
$ws.Range("L117").Value = 85715508
$ws.Range("M117").Value = 2146
$ws.Range("N117").Value = -85722392

$ws.Range("H131").Value = 810.7
$ws.Range("J131").Value = 813.9796
$ws.Range("L131").Value = 2441.9388
$ws.Range("N131").Value = -12521.9388

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3999.3333
$ws.Range("I80").Value = 3999.3333
$ws.Range("K80").Value = 3999.3333
$ws.Range("M80").Value = -3001.3333

$ws.Range("H83").Value = 3999.3333
$ws.Range("I83").Value = 3999.3333
$ws.Range("K83").Value = 19996.6665
$ws.Range("M83").Value = -15004.6665

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 6042.7144
$ws.Range("I32").Value = 4511.8
$ws.Range("J32").Value = 9870
$ws.Range("K32").Value = 4511.8
$ws.Range("L32").Value = 9870
$ws.Range("M32").Value = -4194.8
$ws.Range("N32").Value = -10504

$ws.Range("H132").Value = 2025.0646
$ws.Range("I132").Value = 2078.4
$ws.Range("J132").Value = 1999.6666
$ws.Range("K132").Value = 6235.200000000001
$ws.Range("L132").Value = 5998.9998
$ws.Range("M132").Value = -3705.200000000001
$ws.Range("N132").Value = -11058.9998

$ws.Range("H136").Value = 2325.3215
$ws.Range("I136").Value = 1409.9524
$ws.Range("K136").Value = 4229.857199999999
$ws.Range("M136").Value = -1679.857199999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 14161.667
$ws.Range("I29").Value = 2495
$ws.Range("J29").Value = 19995
$ws.Range("K29").Value = 2495
$ws.Range("L29").Value = 19995
$ws.Range("M29").Value = -2205
$ws.Range("N29").Value = -20575

$ws.Range("H81").Value = 1457.7
$ws.Range("I81").Value = 1508.5555
$ws.Range("K81").Value = 3017.111
$ws.Range("M81").Value = -1956.111

$ws.Range("H84").Value = 1457.7
$ws.Range("I84").Value = 1508.5555
$ws.Range("K84").Value = 15085.555
$ws.Range("M84").Value = -9781.555

$ws.Range("H92").Value = 24850
$ws.Range("J92").Value = 24775
$ws.Range("L92").Value = 24775
$ws.Range("N92").Value = -29767

$ws.Range("H95").Value = 47343.332
$ws.Range("J95").Value = 47343.332
$ws.Range("L95").Value = 47343.332
$ws.Range("N95").Value = -52835.332

$ws.Range("H122").Value = 57178.5
$ws.Range("I122").Value = 98346.875
$ws.Range("K122").Value = 295040.625
$ws.Range("M122").Value = -292590.625

$ws.Range("H132").Value = 1179.878
$ws.Range("I132").Value = 889.3871
$ws.Range("K132").Value = 2668.1613
$ws.Range("M132").Value = -138.1613000000002

$ws.Range("H136").Value = 23150966
$ws.Range("I136").Value = 29242564
$ws.Range("J136").Value = 2890
$ws.Range("K136").Value = 87727692
$ws.Range("L136").Value = 8670
$ws.Range("M136").Value = -87725142
$ws.Range("N136").Value = -13770

